$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row-level fix: row 3 (GRM155R71H104KE14J, 100nF) quantity/subtotal corrected ---
$ws.Cells.Item(3, 4).Value = 16
$ws.Cells.Item(3, 8).Value = 0.496

# --- Delete the duplicate "GRM188R60J106ME47D" / 120nF row (old row 4) entirely. ---
# Everything below shifts up by one row after this.
$ws.Rows.Item(4).Delete()

# --- Row 8 (was row 9 before the delete): GRM1885C1H102JA01D / 1nF price fix ---
$ws.Cells.Item(8, 7).Value = 0.008
$ws.Cells.Item(8, 8).Value = 0.008

# --- Row 11 (was row 12 before the delete): part number + price fix ---
$ws.Cells.Item(11, 3).Value = "GRM188R61E225KA12D"
$ws.Cells.Item(11, 7).Value = 0.102
$ws.Cells.Item(11, 8).Value = 0.204

# --- Row 32 (was row 33 before the delete): fill in missing MPN + pricing ---
$ws.Cells.Item(32, 3).Value = "MSP432E401YTPDT"
$ws.Cells.Item(32, 7).Value = 16.6
$ws.Cells.Item(32, 8).Value = 16.6

# --- Row 40 (was row 41 before the delete): price fix ---
$ws.Cells.Item(40, 7).Value = 0.813
$ws.Cells.Item(40, 8).Value = 0.813

Write-Output "done"
